# Auto-generated update of FFXIV leve profit market-price data
# Applies per-cell value updates (and a few cell clears/additions) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "ALC" = @(
        @{ Cell = "H2"; Value = 356.25 }
        @{ Cell = "I2"; Value = 150 }
        @{ Cell = "J2"; Value = 480 }
        @{ Cell = "K2"; Value = 150 }
        @{ Cell = "L2"; Value = 480 }
        @{ Cell = "M2"; Value = -37 }
        @{ Cell = "N2"; Value = -706 }
        @{ Cell = "H42"; Value = 1666.8334 }
        @{ Cell = "I42"; Value = 1291.8 }
        @{ Cell = "J42"; Value = 2135.625 }
        @{ Cell = "K42"; Value = 3875.4 }
        @{ Cell = "L42"; Value = 6406.875 }
        @{ Cell = "M42"; Value = -3645.4 }
        @{ Cell = "N42"; Value = -6866.875 }
        @{ Cell = "H43"; Value = 3632.8 }
        @{ Cell = "I43"; Value = 3499.6667 }
        @{ Cell = "J43"; Value = 3832.5 }
        @{ Cell = "K43"; Value = 3499.6667 }
        @{ Cell = "L43"; Value = 3832.5 }
        @{ Cell = "M43"; Value = -3430.6667 }
        @{ Cell = "N43"; Value = -3970.5 }
        @{ Cell = "H48"; Value = 0 }
        @{ Cell = "I48"; Value = 0 }
        @{ Cell = "K48"; Value = 0 }
        @{ Cell = "M48"; Value = $null }
        @{ Cell = "H51"; Value = 6989 }
        @{ Cell = "J51"; Value = 6483.5 }
        @{ Cell = "L51"; Value = 6483.5 }
        @{ Cell = "N51"; Value = -7451.5 }
        @{ Cell = "H56"; Value = 0 }
        @{ Cell = "I56"; Value = 0 }
        @{ Cell = "K56"; Value = 0 }
        @{ Cell = "M56"; Value = $null }
        @{ Cell = "H75"; Value = 110473 }
        @{ Cell = "J75"; Value = 121567.6 }
        @{ Cell = "L75"; Value = 121567.6 }
        @{ Cell = "N75"; Value = -123439.6 }
        @{ Cell = "H78"; Value = 110473 }
        @{ Cell = "J78"; Value = 121567.6 }
        @{ Cell = "L78"; Value = 364702.8 }
        @{ Cell = "N78"; Value = -374062.8 }
        @{ Cell = "H103"; Value = 594.4286 }
        @{ Cell = "I103"; Value = 809 }
        @{ Cell = "J103"; Value = 433.5 }
        @{ Cell = "K103"; Value = 2427 }
        @{ Cell = "L103"; Value = 1300.5 }
        @{ Cell = "M103"; Value = -1841 }
        @{ Cell = "N103"; Value = -2472.5 }
        @{ Cell = "H138"; Value = 2355.6814 }
        @{ Cell = "I138"; Value = 2127.9333 }
        @{ Cell = "J138"; Value = 2467.6885 }
        @{ Cell = "K138"; Value = 6383.7999 }
        @{ Cell = "L138"; Value = 7403.065500000001 }
        @{ Cell = "M138"; Value = -1243.7999 }
        @{ Cell = "N138"; Value = -17683.0655 }
        @{ Cell = "H141"; Value = 3353.1875 }
        @{ Cell = "I141"; Value = 3353.1875 }
        @{ Cell = "J141"; Value = 0 }
        @{ Cell = "K141"; Value = 10059.5625 }
        @{ Cell = "L141"; Value = 0 }
        @{ Cell = "M141"; Value = -4879.5625 }
        @{ Cell = "N141"; Value = $null }
    )
    "ARM" = @(
        @{ Cell = "H32"; Value = 2754.9844 }
        @{ Cell = "I32"; Value = 1714.7966 }
        @{ Cell = "K32"; Value = 1714.7966 }
        @{ Cell = "M32"; Value = -1427.7966 }
        @{ Cell = "H61"; Value = 57097.523 }
        @{ Cell = "I61"; Value = 1278.4667 }
        @{ Cell = "K61"; Value = 1278.4667 }
        @{ Cell = "M61"; Value = -1066.4667 }
        @{ Cell = "H74"; Value = 10579.667 }
        @{ Cell = "I74"; Value = 1434 }
        @{ Cell = "J74"; Value = 52878.375 }
        @{ Cell = "K74"; Value = 1434 }
        @{ Cell = "L74"; Value = 52878.375 }
        @{ Cell = "M74"; Value = -560 }
        @{ Cell = "N74"; Value = -54626.375 }
        @{ Cell = "H77"; Value = 10579.667 }
        @{ Cell = "I77"; Value = 1434 }
        @{ Cell = "J77"; Value = 52878.375 }
        @{ Cell = "K77"; Value = 7170 }
        @{ Cell = "L77"; Value = 264391.875 }
        @{ Cell = "M77"; Value = -2802 }
        @{ Cell = "N77"; Value = -273127.875 }
        @{ Cell = "H136"; Value = 57097.523 }
        @{ Cell = "I136"; Value = 1278.4667 }
        @{ Cell = "K136"; Value = 3835.4001 }
        @{ Cell = "M136"; Value = -1285.4001 }
    )
    "BSM" = @(
        @{ Cell = "H99"; Value = 31509.8 }
        @{ Cell = "J99"; Value = 3972.25 }
        @{ Cell = "L99"; Value = 3972.25 }
        @{ Cell = "N99"; Value = -6968.25 }
        @{ Cell = "H132"; Value = 69999 }
        @{ Cell = "J132"; Value = 69999 }
        @{ Cell = "L132"; Value = 69999 }
        @{ Cell = "N132"; Value = -80119 }
    )
    "CRP" = @(
        @{ Cell = "H48"; Value = 48842 }
        @{ Cell = "J48"; Value = 48842 }
        @{ Cell = "L48"; Value = 48842 }
        @{ Cell = "N48"; Value = -49794 }
        @{ Cell = "H99"; Value = 6876645.5 }
        @{ Cell = "I99"; Value = 7124360.5 }
        @{ Cell = "J99"; Value = 6670216.5 }
        @{ Cell = "K99"; Value = 7124360.5 }
        @{ Cell = "L99"; Value = 6670216.5 }
        @{ Cell = "M99"; Value = -7122862.5 }
        @{ Cell = "N99"; Value = -6673212.5 }
        @{ Cell = "H126"; Value = 6876645.5 }
        @{ Cell = "I126"; Value = 7124360.5 }
        @{ Cell = "J126"; Value = 6670216.5 }
        @{ Cell = "K126"; Value = 21373081.5 }
        @{ Cell = "L126"; Value = 20010649.5 }
        @{ Cell = "M126"; Value = -21370611.5 }
        @{ Cell = "N126"; Value = -20015589.5 }
        @{ Cell = "H134"; Value = 41674748 }
        @{ Cell = "I134"; Value = 1881.5 }
        @{ Cell = "K134"; Value = 5644.5 }
        @{ Cell = "M134"; Value = -3109.5 }
    )
    "CUL" = @(
        @{ Cell = "H40"; Value = 339.52 }
        @{ Cell = "I40"; Value = 466.6 }
        @{ Cell = "J40"; Value = 148.9 }
        @{ Cell = "K40"; Value = 1866.4 }
        @{ Cell = "L40"; Value = 595.6 }
        @{ Cell = "M40"; Value = -1797.4 }
        @{ Cell = "N40"; Value = -733.6 }
        @{ Cell = "H68"; Value = 2841.3333 }
        @{ Cell = "J68"; Value = 2841.3333 }
        @{ Cell = "L68"; Value = 8523.999899999999 }
        @{ Cell = "N68"; Value = -10145.9999 }
        @{ Cell = "H71"; Value = 2841.3333 }
        @{ Cell = "J71"; Value = 2841.3333 }
        @{ Cell = "L71"; Value = 25571.9997 }
        @{ Cell = "N71"; Value = -33683.9997 }
        @{ Cell = "H113"; Value = 1270.1111 }
        @{ Cell = "J113"; Value = 1333 }
        @{ Cell = "L113"; Value = 3999 }
        @{ Cell = "N113"; Value = -8339 }
        @{ Cell = "H122"; Value = 7688187.5 }
        @{ Cell = "I122"; Value = 20763454 }
        @{ Cell = "J122"; Value = 1494640.1 }
        @{ Cell = "K122"; Value = 186871086 }
        @{ Cell = "L122"; Value = 13451760.9 }
        @{ Cell = "M122"; Value = -186868636 }
        @{ Cell = "N122"; Value = -13456660.9 }
        @{ Cell = "H129"; Value = 15153210 }
        @{ Cell = "I129"; Value = 2292.75 }
        @{ Cell = "J129"; Value = 45455044 }
        @{ Cell = "K129"; Value = 6878.25 }
        @{ Cell = "L129"; Value = 136365132 }
        @{ Cell = "M129"; Value = -1878.25 }
        @{ Cell = "N129"; Value = -136375132 }
    )
    "GSM" = @(
        @{ Cell = "H11"; Value = 504550.03 }
        @{ Cell = "I11"; Value = 591008.8 }
        @{ Cell = "J11"; Value = 14617 }
        @{ Cell = "K11"; Value = 591008.8 }
        @{ Cell = "L11"; Value = 14617 }
        @{ Cell = "M11"; Value = -590869.8 }
        @{ Cell = "N11"; Value = -14895 }
    )
    "LTW" = @(
        @{ Cell = "H7"; Value = 1708191.9 }
        @{ Cell = "I7"; Value = 2983275 }
        @{ Cell = "K7"; Value = 2983275 }
        @{ Cell = "M7"; Value = -2983163 }
        @{ Cell = "H55"; Value = 1924.909 }
        @{ Cell = "I55"; Value = 1546.9166 }
        @{ Cell = "J55"; Value = 2378.5 }
        @{ Cell = "K55"; Value = 1546.9166 }
        @{ Cell = "L55"; Value = 2378.5 }
        @{ Cell = "M55"; Value = -1373.9166 }
        @{ Cell = "N55"; Value = -2724.5 }
        @{ Cell = "H64"; Value = 300000 }
        @{ Cell = "J64"; Value = 300000 }
        @{ Cell = "L64"; Value = 300000 }
        @{ Cell = "N64"; Value = -300450 }
        @{ Cell = "H67"; Value = 300000 }
        @{ Cell = "J67"; Value = 300000 }
        @{ Cell = "L67"; Value = 300000 }
        @{ Cell = "N67"; Value = -301560 }
        @{ Cell = "H100"; Value = 3993.5 }
        @{ Cell = "I100"; Value = 4158 }
        @{ Cell = "K100"; Value = 4158 }
        @{ Cell = "M100"; Value = -3617 }
        @{ Cell = "H126"; Value = 1708191.9 }
        @{ Cell = "I126"; Value = 2983275 }
        @{ Cell = "K126"; Value = 8949825 }
        @{ Cell = "M126"; Value = -8947355 }
    )
    "WVR" = @(
        @{ Cell = "H68"; Value = 47500 }
        @{ Cell = "J68"; Value = 47500 }
        @{ Cell = "L68"; Value = 47500 }
        @{ Cell = "N68"; Value = -49122 }
        @{ Cell = "H69"; Value = 37687.5 }
        @{ Cell = "J69"; Value = 37687.5 }
        @{ Cell = "L69"; Value = 37687.5 }
        @{ Cell = "N69"; Value = -39185.5 }
        @{ Cell = "H71"; Value = 47500 }
        @{ Cell = "J71"; Value = 47500 }
        @{ Cell = "L71"; Value = 142500 }
        @{ Cell = "N71"; Value = -150612 }
        @{ Cell = "H72"; Value = 37687.5 }
        @{ Cell = "J72"; Value = 37687.5 }
        @{ Cell = "L72"; Value = 113062.5 }
        @{ Cell = "N72"; Value = -120550.5 }
        @{ Cell = "H100"; Value = 771.9 }
        @{ Cell = "I100"; Value = 694.1667 }
        @{ Cell = "J100"; Value = 888.5 }
        @{ Cell = "K100"; Value = 1388.3334 }
        @{ Cell = "L100"; Value = 1777 }
        @{ Cell = "M100"; Value = -847.3334 }
        @{ Cell = "N100"; Value = -2859 }
        @{ Cell = "H132"; Value = 18580.75 }
        @{ Cell = "I132"; Value = 7670.0415 }
        @{ Cell = "K132"; Value = 23010.1245 }
        @{ Cell = "M132"; Value = -20480.1245 }
        @{ Cell = "H133"; Value = 56914.285 }
        @{ Cell = "J133"; Value = 56914.285 }
        @{ Cell = "L133"; Value = 56914.285 }
        @{ Cell = "N133"; Value = -67034.285 }
        @{ Cell = "H135"; Value = 100715 }
        @{ Cell = "J135"; Value = 100715 }
        @{ Cell = "L135"; Value = 100715 }
        @{ Cell = "N135"; Value = -110855 }
    )
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($update in $sheetUpdates[$sheetName]) {
        if ($null -eq $update.Value) {
            $ws.Range($update.Cell).ClearContents()
        } else {
            $ws.Range($update.Cell).Value = $update.Value
        }
    }
}

Write-Host "Updated $($sheetUpdates.Keys.Count) sheets."